# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46009

$ws.Range("B2").Value = 109.35
$ws.Range("C2").Value = 99.61
$ws.Range("D2").Value = 97.34
$ws.Range("E2").Value = 87.78
$ws.Range("F2").Value = 84.25
$ws.Range("G2").Value = 84.79000000000001
$ws.Range("H2").Value = 93.48
$ws.Range("I2").Value = 98.45999999999999
$ws.Range("J2").Value = 105.38
$ws.Range("K2").Value = 92.03
$ws.Range("L2").Value = 76.23
$ws.Range("M2").Value = 46.18
$ws.Range("N2").Value = 40.66
$ws.Range("O2").Value = 30.47
$ws.Range("P2").Value = 32.44
$ws.Range("Q2").Value = 60.98
$ws.Range("R2").Value = 79.79000000000001
$ws.Range("S2").Value = 99.18000000000001
$ws.Range("T2").Value = 103.6
$ws.Range("U2").Value = 108.61
$ws.Range("V2").Value = 110.85
$ws.Range("W2").Value = 109.19
$ws.Range("X2").Value = 98.33
$ws.Range("Y2").Value = 92.23
$ws.Range("Z2").Value = 85.05

$ws.Range("AB2").Value = 102.65
$ws.Range("AD2").Value = 110.02
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 106.1
$ws.Range("AG2").Value = "4h-16h"
